$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: split the old combined "wage period from...to..." label into three cells ---
$ws.Range("A5").Value = "Wage period"
$ws.Range("B5").Value = "From:"
$ws.Range("C5").Value = "To:"

# --- Row 6: rename the "name" header to "Name" ---
$ws.Range("B6").Value = "Name"

# --- Column widths (values chosen so the saved OOXML width lands as close as possible
#     to the target width given this runtime's column-width rounding behaviour) ---
$ws.Columns.Item(1).ColumnWidth = 22.833333333333332
$ws.Columns.Item(2).ColumnWidth = 32.333333333333336
$ws.Columns.Item(3).ColumnWidth = 36.333333333333336
$ws.Columns.Item(4).ColumnWidth = 21.833333333333332
$ws.Columns.Item(9).ColumnWidth = 28.0
$ws.Columns.Item(10).ColumnWidth = 24.333333333333332
$ws.Columns.Item(11).ColumnWidth = 25.666666666666668
$ws.Columns.Item(12).ColumnWidth = 22.333333333333332
$ws.Columns.Item(13).ColumnWidth = 11.666666666666666
$ws.Columns.Item(14).ColumnWidth = 12.333333333333334
$ws.Columns.Item(15).ColumnWidth = 12.833333333333334
$ws.Columns.Item(16).ColumnWidth = 22.333333333333332
$ws.Columns.Item(17).ColumnWidth = 18.0
$ws.Columns.Item(18).ColumnWidth = 22.666666666666668

# --- Sheet view: scroll so column I is the left-most visible column, and select B12 ---
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("B12").Select()
